$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New dSF (column F) values re-pulled from source data.
$newValues = @{
    3  = -1
    4  = 1
    5  = -4
    6  = 2
    7  = -6
    8  = 3
    9  = 9
    10 = -3
    11 = -2
    12 = 1
    13 = 4
    14 = -2
    15 = 2
    16 = 2
    17 = -1
    18 = 5
    19 = 8
    20 = 2
    21 = -4
    22 = 6
    24 = -4
    25 = 3
    26 = 1
    28 = 3
    29 = 3
    30 = 2
    31 = 1
    32 = -1
    33 = -1
    34 = 1
    35 = -5
    36 = 2
    37 = 2
    38 = -2
}

foreach ($rowNum in $newValues.Keys) {
    $ws.Range("F$rowNum").Value = $newValues[$rowNum]
}
